$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ARMS")

# Copy formatting from row 3 into row 4, then set the new row's values
$ws.Range("A3:L3").Copy()
$ws.Range("A4:L4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A4").Value = 42991
$ws.Range("B4").Value = "ARMS3"
$ws.Range("C4").Value = "Test Arms Assessment"
$ws.Range("D4").Value = "Tom Swann"
$ws.Range("E4").Value = 1001
$ws.Range("F4").Value = "C"
$ws.Range("G4").Value = "WMT"
$ws.Range("H4").Value = "A1"
$ws.Range("I4").Value = "C1234569"
$ws.Range("J4").Value = 42991
$ws.Range("K4").Value = "Community"
$ws.Range("L4").Value = 42991

# Update B3: ARMS1 -> ARMS2 (added to shared strings after row 4's new values)
$ws.Range("B3").Value = "ARMS2"

$ws.Range("B3").Select()
